$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. "7.67") are preserved exactly as text instead of being
# auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.289.64'
$ws.Range('E2').Value = '  +4.14%  '
$ws.Range('D3').Value = '3.488.98'
$ws.Range('E3').Value = '  +3.72%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '586.04'
$ws.Range('E5').Value = '  +3.10%  '
$ws.Range('D6').Value = '147.91'
$ws.Range('E8').Value = '  +1.53%  '
$ws.Range('D9').Value = '7.67'
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('D10').Value = '0.126'
$ws.Range('E10').Value = '  +4.80%  '
$ws.Range('D11').Value = '0.399'
$ws.Range('E11').Value = '  +5.00%  '
$ws.Range('D12').Value = '4.085.72'
$ws.Range('D13').Value = '29.82'
$ws.Range('E13').Value = '  +7.72%  '
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').Value = '3.477.75'
$ws.Range('E15').Value = '  +3.56%  '
$ws.Range('E16').Value = '  +4.32%  '
$ws.Range('D17').Value = '63.354.49'
$ws.Range('E17').Value = '  +4.00%  '
$ws.Range('E18').Value = '  +3.57%  '
$ws.Range('E19').Value = '  +6.22%  '
$ws.Range('D20').Value = '9.47'
$ws.Range('E20').Value = '  +6.64%  '
$ws.Range('D21').Value = '393.16'
$ws.Range('E21').Value = '  +3.17%  '
$ws.Range('E22').Value = '  +3.32%  '
$ws.Range('D23').Value = '75.42'
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('E25').Value = '  +9.72%  '
$ws.Range('D26').Value = '3.634.72'
$ws.Range('E26').Value = '  +3.70%  '
$ws.Range('D27').Value = '0.186'
$ws.Range('E27').Value = '  -2.18%  '
$ws.Range('D28').Value = '7.86'
$ws.Range('E28').Value = '  +10.63%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').Value = '8.26'
$ws.Range('E30').Value = '  +5.75%  '
$ws.Range('E31').Value = '  +2.94%  '
$ws.Range('D32').Value = '1.41'
$ws.Range('E32').Value = '  +6.30%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('D34').Value = '23.84'
$ws.Range('E34').Value = '  +4.13%  '
$ws.Range('D35').Value = '32.59'
$ws.Range('E35').Value = '  +28.78%  '
$ws.Range('D36').Value = '5.35'
$ws.Range('E36').Value = '  +9.40%  '
$ws.Range('E37').Value = '  +5.41%  '
$ws.Range('D38').Value = '172.06'
$ws.Range('E38').Value = '  +2.79%  '
$ws.Range('E39').Value = '  +10.17%  '
$ws.Range('D40').Value = '3.527.34'
$ws.Range('E40').Value = '  +3.80%  '
$ws.Range('D41').Value = '0.0770'
$ws.Range('E41').Value = '  +2.17%  '
$ws.Range('E42').Value = '  +4.09%  '
$ws.Range('E43').Value = '  +8.21%  '
$ws.Range('E44').Value = '  +4.39%  '
$ws.Range('D45').Value = '42.55'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('E46').Value = '  +10.91%  '
$ws.Range('D47').Value = '2.614.92'
$ws.Range('E47').Value = '  +6.45%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '23.85'
$ws.Range('E48').Value = '  +8.13%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').Value = '2.31'
$ws.Range('E49').Value = '  +17.46%  '
$ws.Range('D50').Value = '6.77'
$ws.Range('E50').Value = '  +2.54%  '
$ws.Range('E51').Value = '  +5.63%  '

# Restore default (unstyled) cell style for column D now that the
# text values have been written, matching the original workbook
# which had no explicit style index on these cells.
$ws.Range("D2:D51").Style = "Normal"

